# Updating barcode_offset and umi_offset
#
# 1. barcode_offset list: "0" -> "0,38,76", and a new "10,48,86" entry
# 2. umi_offset list: new "1" entry inserted before "Not applicable"
# 3. The two dependent data-validation ranges on the main ATACseq sheet
#    need to grow to cover the newly added list entries.
# 4. .metadata!pav:createdOn timestamp bumped to reflect the re-export.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) barcode_offset sheet: update existing "0" entry, add "10,48,86"
# ---------------------------------------------------------------------
$wsBarcodeOffset = $wb.Worksheets.Item("barcode_offset")

$wsBarcodeOffset.Range("A1").NumberFormat = "@"
$wsBarcodeOffset.Range("A1").Value = "0,38,76"
$wsBarcodeOffset.Range("A1").ClearFormats()

$wsBarcodeOffset.Range("A5").NumberFormat = "@"
$wsBarcodeOffset.Range("A5").Value = "10,48,86"
$wsBarcodeOffset.Range("A5").ClearFormats()

# ---------------------------------------------------------------------
# 2) umi_offset sheet: insert "1" as the new first entry
# ---------------------------------------------------------------------
$wsUmiOffset = $wb.Worksheets.Item("umi_offset")

$wsUmiOffset.Range("A1").Insert()
$wsUmiOffset.Range("A1").NumberFormat = "@"
$wsUmiOffset.Range("A1").Value = "1"
$wsUmiOffset.Range("A1").ClearFormats()

# ---------------------------------------------------------------------
# 3) Grow the data validation ranges that feed off the two lists above.
#    Touching any Validation object rewrites the whole dataValidations
#    collection and resets allowBlank/showErrorMessage to defaults, so
#    re-apply those flags on every validation on the sheet to preserve
#    the original behaviour.
# ---------------------------------------------------------------------
$wsMain = $wb.Worksheets.Item("ATACseq")

$validatedRanges = @(
    "D2:D1001", "E2:E1001", "F2:F1001", "G2:G1001", "H2:H1001", "I2:I1001",
    "J2:J1001", "K2:K1001", "L2:L1001", "O2:O1001", "P2:P1001", "Q2:Q1001",
    "R2:R1001", "S2:S1001", "T2:T1001", "U2:U1001", "V2:V1001", "X2:X1001",
    "Y2:Y1001", "Z2:Z1001", "AA2:AA1001", "AB2:AB1001", "AC2:AC1001",
    "AD2:AD1001", "AE2:AE1001", "AF2:AF1001", "AG2:AG1001", "AH2:AH1001",
    "AJ2:AJ1001", "AK2:AK1001", "AL2:AL1001", "AN2:AN1001", "AO2:AO1001",
    "AR2:AR1001", "AS2:AS1001", "AT2:AT1001"
)

foreach ($sq in $validatedRanges) {
    $validation = $wsMain.Range($sq).Validation
    $validation.IgnoreBlank = $true
    $validation.ShowError = $true
}

$wsMain.Range("O2:O1001").Validation.Formula1 = "'barcode_offset'!`$A`$1:`$A`$5"
$wsMain.Range("R2:R1001").Validation.Formula1 = "'umi_offset'!`$A`$1:`$A`$3"

# ---------------------------------------------------------------------
# 4) Bump the .metadata!pav:createdOn timestamp
# ---------------------------------------------------------------------
$wsMetadata = $wb.Worksheets.Item(".metadata")
$wsMetadata.Range("C2").Value = "2023-10-31T13:53:10-07:00"
